$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C1").Value = "Danh sách sinh viên"
$ws.Range("C6").Select()
